$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 4808.3335
Write-Host "done"
